$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.745.04'
$ws.Range('E2').Value = '  -2.89%  '
$ws.Range('D3').Value = '1.741.89'
$ws.Range('E3').Value = '  -5.48%  '
$ws.Range('D4').Value = '0.9974'
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '237.90'
$ws.Range('E5').Value = '  -9.38%  '
$ws.Range('D6').Value = '0.9990'
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D7').Value = '0.4927'
$ws.Range('E7').Value = '  -7.30%  '
$ws.Range('D8').Value = '41.50'
$ws.Range('E8').Value = '  -7.65%  '
$ws.Range('D9').Value = '0.2661'
$ws.Range('E9').Value = '  -13.42%  '
$ws.Range('D10').Value = '0.06114'
$ws.Range('E10').Value = '  -11.38%  '
$ws.Range('D11').Value = '1.738.96'
$ws.Range('E11').Value = '  -6.05%  '
$ws.Range('D12').Value = '0.06881'
$ws.Range('E12').Value = '  -11.88%  '
$ws.Range('D13').Value = '15.26'
$ws.Range('E13').Value = '  -17.19%  '
$ws.Range('D14').Value = '4.467'
$ws.Range('E14').Value = '  -11.00%  '
$ws.Range('D15').Value = '76.30'
$ws.Range('E15').Value = '  -14.88%  '
$ws.Range('D16').Value = '0.5793'
$ws.Range('E16').Value = '  -23.46%  '
$ws.Range('D17').Value = '0.9962'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').Value = '0.9997'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').Value = '25.770.33'
$ws.Range('E19').Value = '  -2.88%  '
$ws.Range('D20').Value = '11.58'
$ws.Range('E20').Value = '  -17.31%  '
$ws.Range('D21').Value = '0.000006666'
$ws.Range('E21').Value = '  -16.10%  '
$ws.Range('D22').Value = '1.957.59'
$ws.Range('E22').Value = '  -5.73%  '
$ws.Range('D23').Value = '4.046'
$ws.Range('E23').Value = '  -12.41%  '
$ws.Range('D24').Value = '7.990'
$ws.Range('E24').Value = '  -14.25%  '
$ws.Range('D25').Value = '5.056'
$ws.Range('E25').Value = '  -15.57%  '
$ws.Range('D26').Value = '137.40'
$ws.Range('E26').Value = '  -3.43%  '
$ws.Range('D27').Value = '1.501'
$ws.Range('E27').Value = '  -11.11%  '
$ws.Range('D28').Value = '1.835'
$ws.Range('E28').Value = '  -16.28%  '
$ws.Range('D29').Value = '14.78'
$ws.Range('E29').Value = '  -13.18%  '
$ws.Range('D30').Value = '101.74'
$ws.Range('E30').Value = '  -8.60%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '3.733'
$ws.Range('E31').Value = '  -12.53%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.08018'
$ws.Range('E32').Value = '  -8.88%  '
$ws.Range('D33').Value = '3.472'
$ws.Range('E33').Value = '  -15.05%  '
$ws.Range('D34').Value = '0.04443'
$ws.Range('E34').Value = '  -7.85%  '
$ws.Range('D35').Value = '0.9989'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').Value = '2.628'
$ws.Range('E36').Value = '  -10.44%  '
$ws.Range('D37').Value = '0.9757'
$ws.Range('E37').Value = '  -13.98%  '
$ws.Range('D38').Value = '0.5956'
$ws.Range('E38').Value = '  -18.47%  '
$ws.Range('D39').Value = '2.622'
$ws.Range('E39').Value = '  -15.58%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').Value = '104.54'
$ws.Range('E40').Value = '  -3.27%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '1.911'
$ws.Range('E41').Value = '  -17.56%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '0.9991'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.01515'
$ws.Range('E43').Value = '  -11.84%  '
$ws.Range('D44').Value = '5.155'
$ws.Range('E44').Value = '  -12.24%  '
$ws.Range('D45').Value = '0.3780'
$ws.Range('E45').Value = '  -21.28%  '
$ws.Range('D46').Value = '0.7260'
$ws.Range('E46').Value = '  -19.54%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1099'
$ws.Range('E47').Value = '  -11.43%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.05200'
$ws.Range('E48').Value = '  -10.29%  '
$ws.Range('D49').Value = '30.00'
$ws.Range('E49').Value = '  -13.93%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = '5.880'
$ws.Range('E50').Value = '  -21.36%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '52.13'
$ws.Range('E51').Value = '  -13.49%  '
